# Generate Report for Handoff
# - Flip status from "In Translation" to "Ready for handoff" everywhere it appears
# - Bump the handoff timestamps forward by one minute
# - Widen the "Status" / language columns to fit the new, longer text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps: bump by one minute ---
$wsOverview.Range("G2").Value = "2016-08-12 18:47:51"
$wsDeDe.Range("H2").Value = "2016-08-12 18:47:51"
$wsZhCn.Range("H2").Value = "2016-08-12 18:47:44"

# --- Column widths: widen to fit "Ready for handoff" ---
$wsOverview.Cells.Item(1, 5).EntireColumn.ColumnWidth = 16.333333
$wsOverview.Cells.Item(1, 6).EntireColumn.ColumnWidth = 16.333333
$wsZhCn.Cells.Item(1, 3).EntireColumn.ColumnWidth = 16.333333
$wsDeDe.Cells.Item(1, 3).EntireColumn.ColumnWidth = 16.333333
